# Auto-generated Excel COM-interop script to apply Aegis_Profits.xlsx numeric updates
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Cells.Item(12, 8).Value = 26534.736
$ws.Cells.Item(12, 9).Value = 231.11111
$ws.Cells.Item(12, 11).Value = 231.11111
$ws.Cells.Item(12, 13).Value = -61.11111
# Row 69
$ws.Cells.Item(69, 8).Value = 3746.6667
$ws.Cells.Item(69, 9).Value = 4000
$ws.Cells.Item(69, 10).Value = 3620
$ws.Cells.Item(69, 11).Value = 12000
$ws.Cells.Item(69, 12).Value = 10860
$ws.Cells.Item(69, 13).Value = -11126
$ws.Cells.Item(69, 14).Value = -12608
# Row 72
$ws.Cells.Item(72, 8).Value = 3746.6667
$ws.Cells.Item(72, 9).Value = 4000
$ws.Cells.Item(72, 10).Value = 3620
$ws.Cells.Item(72, 11).Value = 36000
$ws.Cells.Item(72, 12).Value = 32580
$ws.Cells.Item(72, 13).Value = -31632
$ws.Cells.Item(72, 14).Value = -41316
# Row 86
$ws.Cells.Item(86, 8).Value = 1674.1364
$ws.Cells.Item(86, 9).Value = 1716.55
$ws.Cells.Item(86, 10).Value = 1250
$ws.Cells.Item(86, 11).Value = 1716.55
$ws.Cells.Item(86, 12).Value = 1250
$ws.Cells.Item(86, 13).Value = -593.55
$ws.Cells.Item(86, 14).Value = -3496
# Row 89
$ws.Cells.Item(89, 8).Value = 1674.1364
$ws.Cells.Item(89, 9).Value = 1716.55
$ws.Cells.Item(89, 10).Value = 1250
$ws.Cells.Item(89, 11).Value = 8582.75
$ws.Cells.Item(89, 12).Value = 6250
$ws.Cells.Item(89, 13).Value = -2966.75
$ws.Cells.Item(89, 14).Value = -17482
# Row 96
$ws.Cells.Item(96, 8).Value = 2237.25
$ws.Cells.Item(96, 9).Value = 950
$ws.Cells.Item(96, 10).Value = 2666.3333
$ws.Cells.Item(96, 11).Value = 2850
$ws.Cells.Item(96, 12).Value = 7998.999899999999
$ws.Cells.Item(96, 13).Value = -1477
$ws.Cells.Item(96, 14).Value = -10744.9999
# Row 112
$ws.Cells.Item(112, 8).Value = 1785.4117
$ws.Cells.Item(112, 10).Value = 1863.4667
$ws.Cells.Item(112, 12).Value = 5590.4001
$ws.Cells.Item(112, 14).Value = -7806.4001
# Row 132
$ws.Cells.Item(132, 8).Value = 13168948
$ws.Cells.Item(132, 9).Value = 14718059
$ws.Cells.Item(132, 10).Value = 1503
$ws.Cells.Item(132, 11).Value = 44154177
$ws.Cells.Item(132, 12).Value = 4509
$ws.Cells.Item(132, 13).Value = -44151647
$ws.Cells.Item(132, 14).Value = -9569
# Row 138
$ws.Cells.Item(138, 8).Value = 3297.25
$ws.Cells.Item(138, 9).Value = 1724.9474
$ws.Cells.Item(138, 11).Value = 5174.8422
$ws.Cells.Item(138, 13).Value = -34.84220000000005

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 24197.867
$ws.Cells.Item(32, 9).Value = 5394.0547
$ws.Cells.Item(32, 10).Value = 161465.7
$ws.Cells.Item(32, 11).Value = 5394.0547
$ws.Cells.Item(32, 12).Value = 161465.7
$ws.Cells.Item(32, 13).Value = -5107.0547
$ws.Cells.Item(32, 14).Value = -162039.7
# Row 122
$ws.Cells.Item(122, 8).Value = 1506.1666
$ws.Cells.Item(122, 9).Value = 1515.7646
$ws.Cells.Item(122, 11).Value = 4547.293799999999
$ws.Cells.Item(122, 13).Value = -2097.293799999999
# Row 132
$ws.Cells.Item(132, 8).Value = 3603.4167
$ws.Cells.Item(132, 9).Value = 3671.4
$ws.Cells.Item(132, 10).Value = 3448.9092
$ws.Cells.Item(132, 11).Value = 11014.2
$ws.Cells.Item(132, 12).Value = 10346.7276
$ws.Cells.Item(132, 13).Value = -8484.200000000001
$ws.Cells.Item(132, 14).Value = -15406.7276

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 20475.37
$ws.Cells.Item(31, 9).Value = 1374.5555
$ws.Cells.Item(31, 10).Value = 31000.307
$ws.Cells.Item(31, 11).Value = 1374.5555
$ws.Cells.Item(31, 12).Value = 31000.307
$ws.Cells.Item(31, 13).Value = -1079.5555
$ws.Cells.Item(31, 14).Value = -31590.307
# Row 34
$ws.Cells.Item(34, 8).Value = 20475.37
$ws.Cells.Item(34, 9).Value = 1374.5555
$ws.Cells.Item(34, 10).Value = 31000.307
$ws.Cells.Item(34, 11).Value = 1374.5555
$ws.Cells.Item(34, 12).Value = 31000.307
$ws.Cells.Item(34, 13).Value = -1172.5555
$ws.Cells.Item(34, 14).Value = -31404.307
# Row 132
$ws.Cells.Item(132, 8).Value = 100006050
$ws.Cells.Item(132, 9).Value = 125007690
$ws.Cells.Item(132, 10).Value = 71432744
$ws.Cells.Item(132, 11).Value = 375023070
$ws.Cells.Item(132, 12).Value = 214298232
$ws.Cells.Item(132, 13).Value = -375020540
$ws.Cells.Item(132, 14).Value = -214303292

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 16
$ws.Cells.Item(16, 8).Value = 480.4
$ws.Cells.Item(16, 10).Value = 467.33334
$ws.Cells.Item(16, 12).Value = 1402.00002
$ws.Cells.Item(16, 14).Value = -1748.00002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Cells.Item(113, 8).Value = 1659.8
$ws.Cells.Item(113, 9).Value = 1000
$ws.Cells.Item(113, 10).Value = 1824.75
$ws.Cells.Item(113, 11).Value = 1000
$ws.Cells.Item(113, 12).Value = 1824.75
$ws.Cells.Item(113, 13).Value = 1170
$ws.Cells.Item(113, 14).Value = -6164.75
# Row 122
$ws.Cells.Item(122, 8).Value = 958.5769
$ws.Cells.Item(122, 9).Value = 877.63635
$ws.Cells.Item(122, 10).Value = 1017.93335
$ws.Cells.Item(122, 11).Value = 2632.90905
$ws.Cells.Item(122, 12).Value = 3053.80005
$ws.Cells.Item(122, 13).Value = -182.9090500000002
$ws.Cells.Item(122, 14).Value = -7953.80005
# Row 126
$ws.Cells.Item(126, 8).Value = 4660
$ws.Cells.Item(126, 10).Value = 4000
$ws.Cells.Item(126, 12).Value = 12000
$ws.Cells.Item(126, 14).Value = -16940
# Row 132
$ws.Cells.Item(132, 8).Value = 2624.8572
$ws.Cells.Item(132, 9).Value = 2414.15
$ws.Cells.Item(132, 10).Value = 3151.625
$ws.Cells.Item(132, 11).Value = 7242.450000000001
$ws.Cells.Item(132, 12).Value = 9454.875
$ws.Cells.Item(132, 13).Value = -4712.450000000001
$ws.Cells.Item(132, 14).Value = -14514.875

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Cells.Item(61, 8).Value = 1709.85
$ws.Cells.Item(61, 9).Value = 1564.7142
$ws.Cells.Item(61, 10).Value = 2048.5
$ws.Cells.Item(61, 11).Value = 1564.7142
$ws.Cells.Item(61, 12).Value = 2048.5
$ws.Cells.Item(61, 13).Value = -1362.7142
$ws.Cells.Item(61, 14).Value = -2452.5
# Row 93
$ws.Cells.Item(93, 8).Value = 2015.1364
$ws.Cells.Item(93, 9).Value = 2704.8
$ws.Cells.Item(93, 10).Value = 1440.4166
$ws.Cells.Item(93, 11).Value = 2704.8
$ws.Cells.Item(93, 12).Value = 1440.4166
$ws.Cells.Item(93, 13).Value = -1456.8
$ws.Cells.Item(93, 14).Value = -3936.4166
# Row 113
$ws.Cells.Item(113, 8).Value = 1709.85
$ws.Cells.Item(113, 9).Value = 1564.7142
$ws.Cells.Item(113, 10).Value = 2048.5
$ws.Cells.Item(113, 11).Value = 1564.7142
$ws.Cells.Item(113, 12).Value = 2048.5
$ws.Cells.Item(113, 13).Value = 605.2858000000001
$ws.Cells.Item(113, 14).Value = -6388.5
# Row 132
$ws.Cells.Item(132, 8).Value = 3952.4644
$ws.Cells.Item(132, 9).Value = 3939.6365
$ws.Cells.Item(132, 10).Value = 3999.5
$ws.Cells.Item(132, 11).Value = 11818.9095
$ws.Cells.Item(132, 12).Value = 11998.5
$ws.Cells.Item(132, 13).Value = -9288.9095
$ws.Cells.Item(132, 14).Value = -17058.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Cells.Item(113, 8).Value = 758.7
$ws.Cells.Item(113, 9).Value = 560.4
$ws.Cells.Item(113, 10).Value = 957
$ws.Cells.Item(113, 11).Value = 1681.2
$ws.Cells.Item(113, 12).Value = 2871
$ws.Cells.Item(113, 13).Value = 488.8000000000002
$ws.Cells.Item(113, 14).Value = -7211
# Row 116
$ws.Cells.Item(116, 8).Value = 40997.5
$ws.Cells.Item(116, 10).Value = 40997.5
$ws.Cells.Item(116, 12).Value = 40997.5
$ws.Cells.Item(116, 14).Value = -50175.5
# Row 122
$ws.Cells.Item(122, 8).Value = 1634.7646
$ws.Cells.Item(122, 9).Value = 1798.8572
$ws.Cells.Item(122, 10).Value = 1519.9
$ws.Cells.Item(122, 11).Value = 5396.571599999999
$ws.Cells.Item(122, 12).Value = 4559.700000000001
$ws.Cells.Item(122, 13).Value = -2946.571599999999
$ws.Cells.Item(122, 14).Value = -9459.700000000001
# Row 126
$ws.Cells.Item(126, 8).Value = 1231.7646
$ws.Cells.Item(126, 9).Value = 1158.6364
$ws.Cells.Item(126, 10).Value = 1365.8334
$ws.Cells.Item(126, 11).Value = 3475.9092
$ws.Cells.Item(126, 12).Value = 4097.5002
$ws.Cells.Item(126, 13).Value = -1005.9092
$ws.Cells.Item(126, 14).Value = -9037.5002
# Row 132
$ws.Cells.Item(132, 8).Value = 2258.875
$ws.Cells.Item(132, 9).Value = 2531.3572
$ws.Cells.Item(132, 10).Value = 1623.0834
$ws.Cells.Item(132, 11).Value = 7594.071599999999
$ws.Cells.Item(132, 12).Value = 4869.2502
$ws.Cells.Item(132, 13).Value = -5064.071599999999
$ws.Cells.Item(132, 14).Value = -9929.2502
# Row 135
$ws.Cells.Item(135, 8).Value = 36452.223
$ws.Cells.Item(135, 10).Value = 36452.223
$ws.Cells.Item(135, 12).Value = 36452.223
$ws.Cells.Item(135, 14).Value = -46592.223

Write-Host "Applied Aegis_Profits.xlsx updates"